$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 301
$ws.Range("G3").Value = 53
$ws.Range("F4").Value = 220
$ws.Range("F5").Value = 245
$ws.Range("F6").Value = 293
$ws.Range("F7").Value = 7388
$ws.Range("F9").Value = 70
$ws.Range("F10").Value = 3303
$ws.Range("F12").Value = 578
$ws.Range("F13").Value = 588
$ws.Range("F14").Value = 432
$ws.Range("F15").Value = 130
$ws.Range("F16").Value = 25
$ws.Range("F17").Value = 747
$ws.Range("F18").Value = 18
$ws.Range("F20").Value = 204
$ws.Range("F22").Value = 224
$ws.Range("F23").Value = 123
$ws.Range("F24").Value = 365
$ws.Range("F25").Value = 129
$ws.Range("F26").Value = 1071
$ws.Range("F27").Value = 73
$ws.Range("F28").Value = 121
$ws.Range("F29").Value = 2104
$ws.Range("F30").Value = 622
$ws.Range("F31").Value = 25
$ws.Range("F34").Value = 574

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 282
$ws.Range("F4").Value = 307
$ws.Range("F5").Value = 313

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 405

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 405
$ws.Range("F3").Value = 301
$ws.Range("G4").Value = 53
$ws.Range("F5").Value = 220
$ws.Range("F6").Value = 245
$ws.Range("F7").Value = 293
$ws.Range("F8").Value = 7388
$ws.Range("F10").Value = 70
$ws.Range("F11").Value = 282
$ws.Range("F12").Value = 3303
$ws.Range("F14").Value = 578
$ws.Range("F15").Value = 588
$ws.Range("F16").Value = 432
$ws.Range("F18").Value = 130
$ws.Range("F19").Value = 25
$ws.Range("F20").Value = 307
$ws.Range("F21").Value = 313
$ws.Range("F23").Value = 747
$ws.Range("F24").Value = 18
$ws.Range("F26").Value = 204
$ws.Range("F31").Value = 224
$ws.Range("F32").Value = 124
$ws.Range("F33").Value = 365
$ws.Range("F34").Value = 129
$ws.Range("F35").Value = 1071
$ws.Range("F36").Value = 73
$ws.Range("F37").Value = 121
$ws.Range("F38").Value = 2104
$ws.Range("F39").Value = 622
$ws.Range("F40").Value = 25
$ws.Range("F43").Value = 574

